# Apply the "Updated cryptos list" refresh: per-cell text values for the
# Price (D) and Volume(1h) (E) columns, plus the Coin/Link swap for rows 49-50.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume figures are stored as text (e.g. "69.507.15", "  -2.20%  ").
# Cells whose new Price value would otherwise be auto-parsed as a plain number
# (e.g. "1.00", "693.68") are switched to the Text format first so Excel keeps
# the literal string instead of silently converting it to a numeric value.

$ws.Range("D2").Value = '69.507.15'
$ws.Range("E2").Value = '  -2.20%  '
$ws.Range("D3").Value = '3.705.10'
$ws.Range("E3").Value = '  -2.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '693.68'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.34'
$ws.Range("E6").Value = '  -5.06%  '
$ws.Range("D7").Value = '3.704.69'
$ws.Range("E7").Value = '  -2.84%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -4.48%  '
$ws.Range("E10").Value = '  -7.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.41'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  -3.35%  '
$ws.Range("E13").Value = '  -4.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.59'
$ws.Range("E14").Value = '  -6.59%  '
$ws.Range("D15").Value = '4.328.19'
$ws.Range("E15").Value = '  -2.95%  '
$ws.Range("D16").Value = '3.701.53'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").Value = '69.621.83'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.33'
$ws.Range("E19").Value = '  -6.77%  '
$ws.Range("E20").Value = '  -7.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '482.43'
$ws.Range("E21").Value = '  -5.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.03'
$ws.Range("E22").Value = '  -6.29%  '
$ws.Range("E23").Value = '  -6.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.15'
$ws.Range("E24").Value = '  -4.64%  '
$ws.Range("E25").Value = '  -8.74%  '
$ws.Range("D26").Value = '3.848.92'
$ws.Range("E26").Value = '  -3.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.47'
$ws.Range("E28").Value = '  -4.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.56'
$ws.Range("E29").Value = '  -8.19%  '
$ws.Range("E30").Value = '  -9.45%  '
$ws.Range("E31").Value = '  -9.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.87'
$ws.Range("E32").Value = '  -7.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.08'
$ws.Range("E33").Value = '  -7.10%  '
$ws.Range("E34").Value = '  -6.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.167'
$ws.Range("E35").Value = '  -4.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = '3.669.07'
$ws.Range("E37").Value = '  -2.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.51'
$ws.Range("E38").Value = '  -7.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.43'
$ws.Range("E39").Value = '  +7.27%  '
$ws.Range("E40").Value = '  -2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0935'
$ws.Range("E41").Value = '  -7.46%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.954'
$ws.Range("E44").Value = '  -6.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '163.60'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.06'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.23'
$ws.Range("E47").Value = '  +2.71%  '
$ws.Range("E48").Value = '  -13.91%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.16'
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("E51").Value = '  -7.92%  '
